# implementazione descaling lato GPU TS, passo a linux per il build
#
# Adds the new "Descaling" label/translation rows to the "Cleaning" sheet
# (label name in column A, English text in column B, Italian text in
# column C for the header row of the new block), and updates the
# active-sheet / selection bookkeeping to match the authored commit:
#   - "Cleaning" (9th tab) becomes the active/selected sheet
#   - "uso comune" (1st tab) is no longer the selected sheet

$wb = $excel.ActiveWorkbook

$wsCleaning = $wb.Worksheets.Item("Cleaning")
$wsUsoComune = $wb.Worksheets.Item("uso comune")

# --- New descaling rows -----------------------------------------------

$wsCleaning.Range("A79").Value = "`$LAB_DESCALING_START"
$wsCleaning.Range("B79").Value = "START DESCALING"
$wsCleaning.Range("C79").Value = "DECALCIFICAZIONE"

# Row 80 intentionally left blank (matches the source workbook layout).
#
# NOTE on write order: below, column B is (re-)used/written before column A
# for rows 81 and onward. This mirrors the original authoring order and
# reproduces the exact shared-string insertion order/index numbering seen
# in the target workbook (new unique strings are appended to the shared
# string table in the order the cells are written).

$wsCleaning.Range("B81").Value = "Press CONTINUE when done"
$wsCleaning.Range("A81").Value = "`$LAB_DESCALING_PRESS_CONTINUE"

$wsCleaning.Range("A82").Value = "`$LAB_DESCALING_OPEN_BOILER_TAP"
$wsCleaning.Range("B82").Value = "Open boiler tap"

$wsCleaning.Range("A83").Value = "`$LAB_DESCALING_CLOSE_BOILER_TAP"
$wsCleaning.Range("B83").Value = "Close boiler tap"

$wsCleaning.Range("B84").Value = "please wait..."
$wsCleaning.Range("A84").Value = "`$LAB_DESCALING_PLS_WAIT"

$wsCleaning.Range("B85").Value = "Emptying hydraulic circuit"
$wsCleaning.Range("A85").Value = "`$LAB_DESCALING_2"

$wsCleaning.Range("B86").Value = "Attach submersible pump to the tank containing descaling detergent"
$wsCleaning.Range("A86").Value = "`$LAB_DESCALING_4"

$wsCleaning.Range("B87").Value = "Filling hydraulic circuit with descaling detergent"
$wsCleaning.Range("A87").Value = "`$LAB_DESCALING_5"

$wsCleaning.Range("B88").Value = "Check descaling detergent level in the air tank"
$wsCleaning.Range("A88").Value = "`$LAB_DESCALING_6"

$wsCleaning.Range("B89").Value = "Starting to fill hydraulic tubes with descaling detergent"
$wsCleaning.Range("A89").Value = "`$LAB_DESCALING_7"

$wsCleaning.Range("B90").Value = "Please wait for the descaling liquid action..."
$wsCleaning.Range("A90").Value = "`$LAB_DESCALING_8"

$wsCleaning.Range("B91").Value = "Descaling liquid starts draining through the nozzles"
$wsCleaning.Range("A91").Value = "`$LAB_DESCALING_9"

$wsCleaning.Range("B92").Value = "Check the liquid colour drained from the nozzles, it defines if descaling process successfully completed.<br>Press CONTINUE button to continue, or press REPEAT to repeat previous steps."
$wsCleaning.Range("A92").Value = "`$LAB_DESCALING_10"

$wsCleaning.Range("B93").Value = "Emptying hydraulic circuit. All descaling liquid drain out through the nozzles"
$wsCleaning.Range("A93").Value = "`$LAB_DESCALING_12"

$wsCleaning.Range("B94").Value = "Change supply source to water tank. Attach submersible pump to the tank containing water"
$wsCleaning.Range("A94").Value = "`$LAB_DESCALING_14"

$wsCleaning.Range("B95").Value = "Hydraulic circuit will be filled with water"
$wsCleaning.Range("A95").Value = "`$LAB_DESCALING_15"

$wsCleaning.Range("B96").Value = "Check water level into the air tank"
$wsCleaning.Range("A96").Value = "`$LAB_DESCALING_16"

$wsCleaning.Range("B97").Value = "Water drained out through nozzles"
$wsCleaning.Range("A97").Value = "`$LAB_DESCALING_17"

$wsCleaning.Range("B98").Value = "Dispense water and test the sample. Place a cup to collect the sample"
$wsCleaning.Range("A98").Value = "`$LAB_DESCALING_18"

$wsCleaning.Range("B99").Value = "Start draining sample through each nozzle"
$wsCleaning.Range("A99").Value = "`$LAB_DESCALING_19"

$wsCleaning.Range("B100").Value = "Check the pH of collected sample. Press CONTINUE button to continue or REPEAT to repeat the previous steps to clean properly the hydraulic circuit."
$wsCleaning.Range("A100").Value = "`$LAB_DESCALING_20"

$wsCleaning.Range("B101").Value = "Descaling procedure finished"
$wsCleaning.Range("A101").Value = "`$LAB_DESCALING_21"

# --- Active sheet / selection bookkeeping ------------------------------

$wsUsoComune.Activate()
$wsUsoComune.Range("B34").Select()

$wsCleaning.Activate()
$wsCleaning.Range("B101").Select()
